$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Header row (row 1): insert "capacity" column after "name", then add metadata columns
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data row (row 2): shift data right to make room for capacity, and append metadata
$ws.Range("B2").Value = "LEXUSES350(客車）"
$ws.Range("C2").Value = 3456
$ws.Range("D2").Value = "溫子苓"
$ws.Range("E2").Value = "100年04月19闩"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 600000
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-02"
$ws.Range("K2").Value = "丁守中"
$ws.Range("L2").Value = 515
$ws.Range("M2").Value = "tmpf49e1"
$ws.Range("N2").Value = 36

$wb.Save()
